$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Trial 24", 0, 0, 0, 0, 0, 136.4511978626251, 0),
    @("Trial 25", 0, 0, 0, 0, 0, 4.844620943069458, 0),
    @("Trial 26", 0, 0, 0, 0, 0, 41.60684609413147, 0),
    @("Trial 27", 0, 0, 0, 0, 0, 0, 7.688782691955566),
    @("Trial 28", 0, 0, 0, 0, 0, 0, 21.1106321811676)
)

$startRow = 25
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
